# "correccion del panel de control"
# Applies the set of text replacements described by the diff.

$d = $word.ActiveDocument

# 1) Contract date near the top
$d.Content.Find.Execute("03/07/2017", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "04/07/2017", 2) | Out-Null

# 2) Surname -> "veira"
$d.Content.Find.Execute("Acosta", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "veira", 2) | Out-Null

# 3) First name -> "beto"
$d.Content.Find.Execute("Alberto", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "beto", 2) | Out-Null

# 4) DNI number (MatchWholeWord avoids touching the longer CUIT number below)
$d.Content.Find.Execute("38259638", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "12528745", 2) | Out-Null

# 5) CUIT number
$d.Content.Find.Execute("2038259638", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "12458748", 2) | Out-Null

# 6) Address / neighborhood
$d.Content.Find.Execute("Barrio Covifol", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Liniers", 2) | Out-Null

# 7) Subject / program name
$d.Content.Find.Execute("Programacion 4", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Programacion 2", 2) | Out-Null

# 8) Service start date (appears twice -> replace all)
$d.Content.Find.Execute("01/07/2017", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "01/08/2017", 2) | Out-Null

# 9) Service end date (appears twice -> replace all)
$d.Content.Find.Execute("10/12/2017", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "15/12/2017", 2) | Out-Null

# 10) Monthly amount
$d.Content.Find.Execute("8500", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "2300", 2) | Out-Null
